# Update NATMI LR-pair data (Ccl20-Cxcr3) with re-run TPM values.
# New data now has a single sending cluster (FAPs) paired with four
# target clusters, including a new "ECs" cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs","Ccl20","Cxcr3","ECs",3,1,1.019213,3.057639,1,1,1,0.5,0.358189,0.716378,0.052978466298774,0.03595390511601538,0.3650708852569999,2.190425311542,0.052978466298774,0.03595390511601538),
    @("FAPs","Ccl20","Cxcr3","Inflammatory-Mac",3,1,1.019213,3.057639,1,1,2,0.6666666666666666,0.574863,1.724589,0.08502595018806307,0.08655445766079338,0.585907842819,5.273170585371,0.08502595018806307,0.08655445766079338),
    @("FAPs","Ccl20","Cxcr3","Neutrophils",3,1,1.019213,3.057639,1,1,3,1,5.590969,16.772907,0.8269404217996371,0.8418062905306279,5.698388287397,51.285494586573,0.8269404217996371,0.8418062905306279),
    @("FAPs","Ccl20","Cxcr3","Resolving-Mac",3,1,1.019213,3.057639,1,1,3,1,0.237009,0.7110270000000001,0.0350551617135259,0.03568534669256324,0.241562653917,2.174063885253,0.0350551617135259,0.03568534669256324)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}
